# Updated cryptos list on Fri Jun 28 21:20:28 UTC 2024 with GitHub Actions
#
# Refreshes the price (column D) and 1h volume/change (column E) figures for
# the crypto ranking sheet. Two pairs of rows also swapped rank order
# (Avalanche/TRON at rows 14-15, and NEARProtocol/ImmutableX at rows 39-40);
# those are written as full row updates (B, C, D, E) below.
#
# A leading apostrophe is used for column D values that would otherwise be
# auto-interpreted as numbers by Excel (e.g. "567.38", "0.123", "1.00"),
# so they stay plain text with their exact original formatting (trailing
# zeros, fixed decimal places, etc.) preserved - matching how this sheet's
# Price column is stored (plain text, not numeric) everywhere else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.273.25'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '3.381.11'
$ws.Range("E3").Value = '  -1.85%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''567.38'
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").Value = '''139.90'
$ws.Range("E6").Value = '  -6.54%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.383.77'
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").Value = '''0.472'
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").Value = '''7.45'
$ws.Range("E10").Value = '  -5.03%  '
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = '''0.388'
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D13").Value = '3.956.06'
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '''0.123'
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '''27.98'
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = '3.383.86'
$ws.Range("D17").Value = '''0.0000169'
$ws.Range("E17").Value = '  -3.42%  '
$ws.Range("D18").Value = '60.323.38'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("D19").Value = '''6.20'
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").Value = '''13.97'
$ws.Range("E20").Value = '  -2.37%  '
$ws.Range("D21").Value = '''9.04'
$ws.Range("E21").Value = '  -4.85%  '
$ws.Range("D22").Value = '''386.51'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = '''0.557'
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("D24").Value = '''73.29'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '''0.0000115'
$ws.Range("E26").Value = '  -6.58%  '
$ws.Range("D27").Value = '3.530.52'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '''0.178'
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("D30").Value = '''7.35'
$ws.Range("E30").Value = '  -5.27%  '
$ws.Range("D31").Value = '''7.93'
$ws.Range("E31").Value = '  -3.86%  '
$ws.Range("D32").Value = '''2.13'
$ws.Range("E32").Value = '  -2.24%  '
$ws.Range("D33").Value = '''1.40'
$ws.Range("E33").Value = '  -7.62%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '''23.60'
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").Value = '3.410.37'
$ws.Range("E36").Value = '  -1.74%  '
$ws.Range("D37").Value = '''6.88'
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").Value = '''168.01'
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '''1.49'
$ws.Range("E39").Value = '  -5.27%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '''4.90'
$ws.Range("E40").Value = '  -8.01%  '
$ws.Range("D41").Value = '''0.0771'
$ws.Range("E41").Value = '  -2.77%  '
$ws.Range("D42").Value = '''26.99'
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").Value = '''0.781'
$ws.Range("E43").Value = '  -1.75%  '
$ws.Range("D44").Value = '''0.999'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '''4.42'
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("D46").Value = '''1.69'
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").Value = '''41.26'
$ws.Range("E47").Value = '  -2.41%  '
$ws.Range("D48").Value = '2.514.36'
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("E49").Value = '  -4.08%  '
$ws.Range("D50").Value = '''6.78'
$ws.Range("E50").Value = '  -4.23%  '
$ws.Range("D51").Value = '''22.91'
$ws.Range("E51").Value = '  -1.62%  '
